$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the bottom of the "feature of interest" table.
# Column A holds a shared gml_id string (re-using existing ones, plus two
# brand-new BID_* values), column B holds a plain numeric id.
$newRows = @(
    @{ Row = 65; A = "BID_6e5c165d-9eb9-4b3d-9a4b-76958e13b4df"; B = 834 },
    @{ Row = 66; A = "BID_04b55dd6-c136-49a6-b142-723c0eb5ee89"; B = 833 },
    @{ Row = 67; A = "BID_1b33d422-98b1-4781-82f2-243f308b4505"; B = 832 },
    @{ Row = 68; A = "BID_059d0a21-a2a7-4aa5-b427-2975517f09ba"; B = 831 },
    @{ Row = 69; A = "BID_1b33d422-98b1-4781-82f2-243f308b4505"; B = 830 },
    @{ Row = 70; A = "BID_d3ce9c37-e52e-44c4-a1fe-5b5e23e23ffb"; B = 829 },
    @{ Row = 71; A = "BID_1b33d422-98b1-4781-82f2-243f308b4505"; B = 828 },
    @{ Row = 72; A = "BID_0a3a6049-77d6-4c9e-8486-c3469fe49cc7"; B = 820 },
    @{ Row = 73; A = "BID_d3ce9c37-e52e-44c4-a1fe-5b5e23e23ffb"; B = 821 },
    @{ Row = 74; A = "BID_c73cd9f0-a016-4905-bed5-a76622d0c010"; B = 822 },
    @{ Row = 75; A = "BID_609039c8-7b22-4349-8f6b-8d0c76d5970a"; B = 823 },
    @{ Row = 76; A = "BID_95c19986-26fe-47ea-8181-469586a18581"; B = 826 }
)

# Pull the formatting used by the rest of the table (row 60) down onto every
# new row first, so the new cells line up with the existing style ("s=2").
$formatSource = $ws.Range("A60:B60")

foreach ($entry in $newRows) {
    $destRow = $ws.Range("A" + $entry.Row + ":B" + $entry.Row)
    $formatSource.Copy()
    $destRow.PasteSpecial(-4122)
    $ws.Cells.Item($entry.Row, 1).Value = $entry.A
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
}

# Row 75's A cell is styled like the header cell (A1, "s=1") rather than the
# regular body style.
$ws.Range("A1").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Cells.Item(75, 1).Value = "BID_609039c8-7b22-4349-8f6b-8d0c76d5970a"

$excel.CutCopyMode = 0
